# Insert a new data row at row 54 (weekly price-report update).
# This pushes the previous rows 54-77 down to 55-78 (unchanged), and
# populates the newly inserted row 54 with this week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(54).Insert()

$ws.Cells.Item(54, 1).Value  = 1
$ws.Cells.Item(54, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(54, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(54, 4).Value  = 44523
$ws.Cells.Item(54, 5).Value  = 15
$ws.Cells.Item(54, 6).Value  = "Fruta"
$ws.Cells.Item(54, 7).Value  = 100102
$ws.Cells.Item(54, 8).Value  = "Cítricos"
$ws.Cells.Item(54, 9).Value  = 100102004
$ws.Cells.Item(54, 10).Value = "Mandarina"
$ws.Cells.Item(54, 11).Value = "Clementina"
$ws.Cells.Item(54, 12).Value = "Segunda"
$ws.Cells.Item(54, 13).Value = 300
$ws.Cells.Item(54, 14).Value = 11000
$ws.Cells.Item(54, 15).Value = 12000
$ws.Cells.Item(54, 16).Value = 11500
$ws.Cells.Item(54, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(54, 18).Value = "Región Metropolitana"
$ws.Cells.Item(54, 19).Value = 575
$ws.Cells.Item(54, 20).Value = 20
